# Repull data, push all data, mean calculation
# Update the dSF ("F") column values for several rows to match the
# re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value  = -1
$ws.Range("F7").Value  = 2
$ws.Range("F8").Value  = 3
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 3
$ws.Range("F26").Value = 0
$ws.Range("F31").Value = -4
